$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A35").Value = "CF"
$ws.Range("B35").Value = "DBI202"
$ws.Range("C35").Value = "DBI202"
$ws.Range("D35").Value = 30
